# Weekly update: insert a new price-report row for the latest week at the
# top of the data block (row 96), pushing the existing historical rows
# down by one. The last existing row ends up duplicated one row lower
# (row 259), exactly mirroring the previous last row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 96; Excel shifts rows 96:258 down to 97:259
# (including their styles), and grows the used range automatically.
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row with this week's figures.
$ws.Cells.Item(96, 1).Value  = 8
$ws.Cells.Item(96, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(96, 3).Value  = "Coquimbo"
$ws.Cells.Item(96, 4).Value  = 44665
$ws.Cells.Item(96, 5).Value  = 4
$ws.Cells.Item(96, 6).Value  = 100112012
$ws.Cells.Item(96, 7).Value  = "Espinaca"
$ws.Cells.Item(96, 8).Value  = "Sin especificar"
$ws.Cells.Item(96, 9).Value  = "Primera"
$ws.Cells.Item(96, 10).Value = 2560
$ws.Cells.Item(96, 11).Value = 450
$ws.Cells.Item(96, 12).Value = 500
$ws.Cells.Item(96, 13).Value = 475
$ws.Cells.Item(96, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(96, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(96, 16).Value = 950
$ws.Cells.Item(96, 17).Value = 0.5
$ws.Cells.Item(96, 18).Value = "Hortaliza"

# Make sure the date column keeps the same date number format as the rest
# of column D (style index "2" in the original file).
$ws.Cells.Item(96, 4).NumberFormat = $ws.Cells.Item(97, 4).NumberFormat
